$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 361
$ws.Range("I9").Value = 292.7143
$ws.Range("K9").Value = 292.7143
$ws.Range("M9").Value = -123.7143
$ws.Range("H12").Value = 790
$ws.Range("I12").Value = 324.57144
$ws.Range("K12").Value = 324.57144
$ws.Range("M12").Value = -154.57144
$ws.Range("H19").Value = 4466.7144
$ws.Range("J19").Value = 4460.2
$ws.Range("L19").Value = 4460.2
$ws.Range("N19").Value = -4810.2
$ws.Range("H28").Value = 670.5833
$ws.Range("I28").Value = 715.875
$ws.Range("J28").Value = 580
$ws.Range("K28").Value = 715.875
$ws.Range("L28").Value = 580
$ws.Range("M28").Value = -230.875
$ws.Range("N28").Value = -1550
$ws.Range("H33").Value = 343.82352
$ws.Range("I33").Value = 137.71428
$ws.Range("K33").Value = 137.71428
$ws.Range("M33").Value = 91.28572
$ws.Range("H41").Value = 414.3158
$ws.Range("I41").Value = 383.5
$ws.Range("K41").Value = 383.5
$ws.Range("M41").Value = 56.5
$ws.Range("H88").Value = 4559.5
$ws.Range("J88").Value = 4416.1665
$ws.Range("L88").Value = 4416.1665
$ws.Range("N88").Value = -5228.1665
$ws.Range("H91").Value = 4559.5
$ws.Range("J91").Value = 4416.1665
$ws.Range("L91").Value = 4416.1665
$ws.Range("N91").Value = -7224.1665
$ws.Range("H118").Value = 620
$ws.Range("I118").Value = 620
$ws.Range("K118").Value = 1860
$ws.Range("M118").Value = -203
$ws.Range("H138").Value = 2554.46
$ws.Range("I138").Value = 2511.4546
$ws.Range("J138").Value = 2566.5898
$ws.Range("K138").Value = 7534.3638
$ws.Range("L138").Value = 7699.769400000001
$ws.Range("M138").Value = -2394.3638
$ws.Range("N138").Value = -17979.7694

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6087.1763
$ws.Range("I32").Value = 5927.469
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 5927.469
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -5640.469
$ws.Range("N32").Value = -10574
$ws.Range("H61").Value = 6024
$ws.Range("I61").Value = 4904.5625
$ws.Range("J61").Value = 11994.333
$ws.Range("K61").Value = 4904.5625
$ws.Range("L61").Value = 11994.333
$ws.Range("M61").Value = -4692.5625
$ws.Range("N61").Value = -12418.333
$ws.Range("H63").Value = 3460.3333
$ws.Range("I63").Value = 2731.5
$ws.Range("K63").Value = 2731.5
$ws.Range("M63").Value = -2045.5
$ws.Range("H66").Value = 3460.3333
$ws.Range("I66").Value = 2731.5
$ws.Range("K66").Value = 13657.5
$ws.Range("M66").Value = -10225.5
$ws.Range("H132").Value = 5861.75
$ws.Range("I132").Value = 4578.45
$ws.Range("J132").Value = 12278.25
$ws.Range("K132").Value = 13735.35
$ws.Range("L132").Value = 36834.75
$ws.Range("M132").Value = -11205.35
$ws.Range("N132").Value = -41894.75
$ws.Range("H136").Value = 6024
$ws.Range("I136").Value = 4904.5625
$ws.Range("J136").Value = 11994.333
$ws.Range("K136").Value = 14713.6875
$ws.Range("L136").Value = 35982.999
$ws.Range("M136").Value = -12163.6875
$ws.Range("N136").Value = -41082.999
$ws.Range("H139").Value = 141958
$ws.Range("J139").Value = 174947.5
$ws.Range("L139").Value = 174947.5
$ws.Range("N139").Value = -185227.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5187.294
$ws.Range("I20").Value = 4218.6
$ws.Range("J20").Value = 6571.143
$ws.Range("K20").Value = 4218.6
$ws.Range("L20").Value = 6571.143
$ws.Range("M20").Value = -3971.6
$ws.Range("N20").Value = -7065.143
$ws.Range("H86").Value = 1908.2858
$ws.Range("I86").Value = 1708.9231
$ws.Range("K86").Value = 1708.9231
$ws.Range("M86").Value = -585.9231
$ws.Range("H89").Value = 1908.2858
$ws.Range("I89").Value = 1708.9231
$ws.Range("K89").Value = 8544.6155
$ws.Range("M89").Value = -2928.6155
$ws.Range("H99").Value = 4788.25
$ws.Range("I99").Value = 3474.3076
$ws.Range("K99").Value = 3474.3076
$ws.Range("M99").Value = -1976.3076
$ws.Range("H140").Value = 86936.5
$ws.Range("J140").Value = 86936.5
$ws.Range("L140").Value = 86936.5
$ws.Range("N140").Value = -97296.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 102.25
$ws.Range("J7").Value = 99.5
$ws.Range("L7").Value = 99.5
$ws.Range("N7").Value = -325.5
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H16").Value = 3010.9312
$ws.Range("I16").Value = 2445.75
$ws.Range("K16").Value = 2445.75
$ws.Range("M16").Value = -2158.75
$ws.Range("H113").Value = 3010.9312
$ws.Range("I113").Value = 2445.75
$ws.Range("K113").Value = 2445.75
$ws.Range("M113").Value = -275.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20814414
$ws.Range("I4").Value = 18834698
$ws.Range("K4").Value = 56504094
$ws.Range("M4").Value = -56503982
$ws.Range("H12").Value = 549.05884
$ws.Range("J12").Value = 530.25
$ws.Range("L12").Value = 1590.75
$ws.Range("N12").Value = -1936.75
$ws.Range("H70").Value = 2222
$ws.Range("I70").Value = 2222
$ws.Range("K70").Value = 6666
$ws.Range("M70").Value = -6351
$ws.Range("H73").Value = 2222
$ws.Range("I73").Value = 2222
$ws.Range("K73").Value = 6666
$ws.Range("M73").Value = -5574
$ws.Range("H76").Value = 2499998
$ws.Range("I76").Value = 2499998
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 7499994
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -7499611
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 2499998
$ws.Range("I79").Value = 2499998
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 7499994
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -7498668
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 125
$ws.Range("I10").Value = 125
$ws.Range("K10").Value = 125
$ws.Range("M10").Value = 44
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H70").Value = 14250.375
$ws.Range("J70").Value = 16500
$ws.Range("L70").Value = 16500
$ws.Range("N70").Value = -17040
$ws.Range("H73").Value = 14250.375
$ws.Range("J73").Value = 16500
$ws.Range("L73").Value = 16500
$ws.Range("N73").Value = -18372
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H126").Value = 3998.889
$ws.Range("I126").Value = 3750
$ws.Range("J126").Value = 4070
$ws.Range("K126").Value = 11250
$ws.Range("L126").Value = 12210
$ws.Range("M126").Value = -8780
$ws.Range("N126").Value = -17150
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H68").Value = 1175
$ws.Range("I68").Value = 1175
$ws.Range("K68").Value = 1175
$ws.Range("M68").Value = -426
$ws.Range("H71").Value = 1175
$ws.Range("I71").Value = 1175
$ws.Range("K71").Value = 5875
$ws.Range("M71").Value = -2131
$ws.Range("H93").Value = 2523.2856
$ws.Range("I93").Value = 2637.4443
$ws.Range("J93").Value = 2317.8
$ws.Range("K93").Value = 2637.4443
$ws.Range("L93").Value = 2317.8
$ws.Range("M93").Value = -1389.4443
$ws.Range("N93").Value = -4813.8
$ws.Range("H100").Value = 5544.8
$ws.Range("I100").Value = 6362
$ws.Range("K100").Value = 6362
$ws.Range("M100").Value = -5821
$ws.Range("H120").Value = 68998
$ws.Range("J120").Value = 68998
$ws.Range("L120").Value = 68998
$ws.Range("N120").Value = -78674
$ws.Range("H132").Value = 10085.814
$ws.Range("I132").Value = 10891.954
$ws.Range("J132").Value = 6538.8
$ws.Range("K132").Value = 32675.862
$ws.Range("L132").Value = 19616.4
$ws.Range("M132").Value = -30145.862
$ws.Range("N132").Value = -24676.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 20002.777
$ws.Range("J15").Value = 7010
$ws.Range("L15").Value = 7010
$ws.Range("N15").Value = -7586
$ws.Range("H20").Value = 9015
$ws.Range("J20").Value = 9015
$ws.Range("L20").Value = 9015
$ws.Range("N20").Value = -9495
$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180
